$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "Price" (column D) cells hold plain text that merely looks numeric
# (e.g. "48.49"). Excel/COM auto-detects such strings on assignment and
# silently coerces the cell to a real Number (which can also lose a
# significant trailing zero, e.g. "10.90" -> 10.9). To keep these cells as
# text - matching the original inline-string cells - force Text format
# before writing the value, then clear the format again so the cell keeps
# no explicit style, same as in the source workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.314.74"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.62"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.83%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4524"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3865"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.49"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -8.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07919"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.46%  "
$ws.Range("E11").Value = "  -3.43%  "
$ws.Range("E12").Value = "  -4.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.856.41"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.901"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.118"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "85.73"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001024"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06548"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.03"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.51%  "
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.502"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.311.05"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.90"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.284"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.069.78"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.59"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.85"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.059"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.442"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.59%  "
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09293"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9310"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.11%  "
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.259"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.223"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02221"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05989"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.100"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -11.36%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5899"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1882"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.07"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -9.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.279"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5609"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.03"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.370"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.912"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -6.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06747"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.63"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.20%  "
